# Apply the "output generated at 456a3b4" update to the 苏州-漫展信息 workbook.
# Sheet "展览" (index 1) and Sheet "全部类型" (index 4) both contain (mostly)
# duplicated event listings that need the same numeric bumps, one picture URL
# swap, and a content swap between two adjacent rows (the "萤火" expo and the
# "白日梦想" ONLY show traded places, with the "萤火" row also picking up a
# refreshed "want to go" count and an updated start date).

$wb = $excel.ActiveWorkbook

# ---------- Sheet "展览" ----------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value = 81
$ws1.Cells.Item(3, 6).Value = 143
$ws1.Cells.Item(4, 6).Value = 2106
$ws1.Cells.Item(5, 6).Value = 381

$ws1.Cells.Item(8, 6).Value = 2098
$ws1.Cells.Item(8, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/mi1bBvpO1715929055334.jpeg"

$ws1.Cells.Item(9, 6).Value = 10885

# Row 14 becomes the "萤火国潮文化节" listing
# (leading apostrophe keeps the cell as plain text instead of an auto-detected
# date, matching how every other date column in this sheet is stored)
$ws1.Cells.Item(14, 2).Value = "'2024-07-19"
$ws1.Cells.Item(14, 3).Value = "苏州·萤火国潮文化节动漫品牌博览会"
$ws1.Cells.Item(14, 4).Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws1.Cells.Item(14, 5).Value = "2024.07.19 10:00-07.21 17:00"
$ws1.Cells.Item(14, 6).Value = 9935
$ws1.Cells.Item(14, 7).Value = 60
$ws1.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83301"
$ws1.Cells.Item(14, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/Eh06dOvF1715926655440.jpeg"

# Row 15 becomes the "白日梦想7.20全职猎人ONLY展" listing
$ws1.Cells.Item(15, 3).Value = "苏州·白日梦想7.20全职猎人ONLY展"
$ws1.Cells.Item(15, 4).Value = "金芳路与新发路交叉口东南120米 万龙大厦"
$ws1.Cells.Item(15, 5).Value = "2024.07.20 09:00-07.20 17:00"
$ws1.Cells.Item(15, 6).Value = 425
$ws1.Cells.Item(15, 7).Value = 72
$ws1.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83508"
$ws1.Cells.Item(15, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg"

$ws1.Cells.Item(16, 6).Value = 1120
$ws1.Cells.Item(18, 6).Value = 5336
$ws1.Cells.Item(20, 6).Value = 3376

# ---------- Sheet "全部类型" ----------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value = 81
$ws4.Cells.Item(3, 6).Value = 143
$ws4.Cells.Item(4, 6).Value = 2106
$ws4.Cells.Item(5, 6).Value = 381

$ws4.Cells.Item(9, 6).Value = 2098
$ws4.Cells.Item(9, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/mi1bBvpO1715929055334.jpeg"

$ws4.Cells.Item(12, 6).Value = 10885

# Row 17 becomes the "萤火国潮文化节" listing
$ws4.Cells.Item(17, 2).Value = "'2024-07-19"
$ws4.Cells.Item(17, 3).Value = "苏州·萤火国潮文化节动漫品牌博览会"
$ws4.Cells.Item(17, 4).Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws4.Cells.Item(17, 5).Value = "2024.07.19 10:00-07.21 17:00"
$ws4.Cells.Item(17, 6).Value = 9935
$ws4.Cells.Item(17, 7).Value = 60
$ws4.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83301"
$ws4.Cells.Item(17, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/Eh06dOvF1715926655440.jpeg"

# Row 18 becomes the "白日梦想7.20全职猎人ONLY展" listing
$ws4.Cells.Item(18, 3).Value = "苏州·白日梦想7.20全职猎人ONLY展"
$ws4.Cells.Item(18, 4).Value = "金芳路与新发路交叉口东南120米 万龙大厦"
$ws4.Cells.Item(18, 5).Value = "2024.07.20 09:00-07.20 17:00"
$ws4.Cells.Item(18, 6).Value = 425
$ws4.Cells.Item(18, 7).Value = 72
$ws4.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83508"
$ws4.Cells.Item(18, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg"

$ws4.Cells.Item(19, 6).Value = 1120
$ws4.Cells.Item(21, 6).Value = 5336
$ws4.Cells.Item(23, 6).Value = 3376
